$wb = $excel.ActiveWorkbook

# Sheet "展览" - update "想去人数" (want-to-go count) values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 223
$ws1.Range("F4").Value = 2540
$ws1.Range("F6").Value = 546

# Sheet "全部类型" - same underlying rows, update matching values
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 223
$ws4.Range("F6").Value = 2540
$ws4.Range("F8").Value = 546
